$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.395.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "'2.016.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'260.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.02%  "
$ws.Range("D6").Value = "'0.620"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'57.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.00%  "
$ws.Range("E9").Value = "  -3.95%  "
$ws.Range("E10").Value = "  -4.36%  "
$ws.Range("E11").Value = "  -3.02%  "
$ws.Range("D12").Value = "'14.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.61%  "
$ws.Range("D13").Value = "'2.312.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").Value = "'21.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.89%  "
$ws.Range("D15").Value = "'0.799"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.83%  "
$ws.Range("E16").Value = "  -5.46%  "
$ws.Range("D17").Value = "'2.035.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "'37.344.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").Value = "'70.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("D20").Value = "'0.0₃0839"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.61%  "
$ws.Range("D21").Value = "'233.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("D22").Value = "'5.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.74%  "
$ws.Range("D23").Value = "'2.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.82%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").Value = "'164.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").Value = "'9.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.31%  "
$ws.Range("D28").Value = "'19.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.44%  "
$ws.Range("E29").Value = "  -5.48%  "
$ws.Range("E30").Value = "  -5.31%  "
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("D32").Value = "'4.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.94%  "
$ws.Range("D33").Value = "'0.0643"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.04%  "
$ws.Range("D34").Value = "'4.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").Value = "'2.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.81%  "
$ws.Range("D36").Value = "'1.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "'3.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.47%  "
$ws.Range("D39").Value = "'5.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "  +3.51%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "'0.0213"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("E43").Value = "  -5.31%  "
$ws.Range("D44").Value = "'1.436.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.28%  "
$ws.Range("E45").Value = "  -8.42%  "
$ws.Range("D46").Value = "'89.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.31%  "
$ws.Range("E47").Value = "  -3.33%  "
$ws.Range("D48").Value = "'2.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("D49").Value = "'7.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.43%  "
$ws.Range("D50").Value = "'2.203.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("E51").Value = "  -10.36%  "
